$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.733.92"
$ws.Range("E2").Value = "  -3.48%  "

$ws.Range("D3").Value = "2.902.23"
$ws.Range("E3").Value = "  -4.25%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.89"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.46"
$ws.Range("E6").Value = "  -6.17%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -2.77%  "

$ws.Range("D9").Value = "2.902.83"
$ws.Range("E9").Value = "  -4.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -2.37%  "

$ws.Range("E11").Value = "  -4.94%  "

$ws.Range("E12").Value = "  -4.07%  "

$ws.Range("E13").Value = "  -3.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.38"
$ws.Range("E14").Value = "  -6.65%  "

$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "3.383.36"
$ws.Range("E16").Value = "  -4.24%  "

$ws.Range("D17").Value = "60.705.00"
$ws.Range("E17").Value = "  -3.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.70"
$ws.Range("E18").Value = "  -5.32%  "

$ws.Range("D19").Value = "2.904.85"
$ws.Range("E19").Value = "  -4.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.80"
$ws.Range("E20").Value = "  -5.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").Value = "  -5.04%  "

$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("E23").Value = "  -5.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.57"
$ws.Range("E24").Value = "  -3.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.87"
$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("E26").Value = "  -2.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.89"
$ws.Range("E27").Value = "  -4.43%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  -3.37%  "

$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.42"
$ws.Range("E33").Value = "  -4.11%  "

$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("D35").Value = "0.0₃0874"
$ws.Range("E35").Value = "  +1.89%  "

$ws.Range("E36").Value = "  -2.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.60"
$ws.Range("E37").Value = "  -5.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  -5.89%  "

$ws.Range("E39").Value = "  -3.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.49"
$ws.Range("E40").Value = "  -2.00%  "

$ws.Range("E41").Value = "  -4.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.60"
$ws.Range("E42").Value = "  -5.82%  "

$ws.Range("E43").Value = "  -3.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.36"
$ws.Range("E44").Value = "  -6.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "377.21"
$ws.Range("E45").Value = "  -3.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("E46").Value = "  -3.14%  "

$ws.Range("D47").Value = "2.691.92"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.38"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.29"
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("E51").Value = "  -2.70%  "
